# edit.ps1
# Applies the commit: split the "20/09/23" date run into three runs reflecting
# the year fix (23 -> 2023), and appends the journal entries for days 2-6
# (21/09 through 25/09/2023), matching the target OOXML diff.

$d = $word.ActiveDocument

# --- Step 1: fix "20/09/23 (jour 1 ) : " paragraph, splitting the date run into 3 runs ---
# (the year "23" becomes "2023" by inserting "20"; original run boundaries for
#  "(jour" / " 1 ) : " are preserved via a full-paragraph OOXML replace so that
#  the resulting runs are NOT re-merged by the editor)
$para3 = $d.Paragraphs.Item(3)
$para3Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>20/09/</w:t></w:r><w:r><w:t>20</w:t></w:r><w:r><w:t xml:space="preserve">23 </w:t></w:r><w:r><w:t>(jour</w:t></w:r><w:r><w:t xml:space="preserve"> 1 ) : </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para3.Range.InsertXML($para3Xml)

# --- Step 2: append the new journal entries (days 2 through 6) at the end of the document ---
$endPos = $d.Content.End
$tailRange = $d.Range($endPos, $endPos)
$tailXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t xml:space="preserve">21/09/2023 (jour 2) : </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Mise en place du </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>form</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> menu et  début de mise en place sur la page du jeux </w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">22/09/2023 (jour 3) : </w:t></w:r></w:p><w:p><w:r><w:t>Ajout du code de base de fonctionnement du jeu (gestion des tir , conditions victoire/défaite , collision)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>23/09.2023 (jour 4) :</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Debug</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> a la recherche des lag (solution </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>backround</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> image responsable des lag) </w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">24/09/2023 (jour 5) : </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Modification du code pour faire apparaître les ennemis sur 2 lignes et début de gestion des retours a la ligne </w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">25/09/2023  (jour 6) : </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Modification du code avec le passage en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>du</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tableau d’</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>enemies</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en 2 l</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> , changement de la gestion du retour a la ligne des ennemis ainsi que gestion de l’apparition des ennemis avec 2 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>list</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">, ajout de la gestion des point de vie pour le jouer et les </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>invaders</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tailRange.InsertXML($tailXml)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
